$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 60692.6
$ws.Range("I17").Value = 381.5
$ws.Range("J17").Value = 69971.234
$ws.Range("K17").Value = 1144.5
$ws.Range("L17").Value = 209913.702
$ws.Range("M17").Value = -976.5
$ws.Range("N17").Value = -210249.702
# Row 40
$ws.Range("H40").Value = 4360
$ws.Range("I40").Value = 2293.3333
$ws.Range("J40").Value = 5600
$ws.Range("K40").Value = 2293.3333
$ws.Range("L40").Value = 5600
$ws.Range("M40").Value = -2118.3333
$ws.Range("N40").Value = -5950
# Row 64
$ws.Range("H64").Value = 6999.6665
$ws.Range("J64").Value = 6999.6665
$ws.Range("L64").Value = 6999.6665
$ws.Range("N64").Value = -7495.6665
# Row 67
$ws.Range("H67").Value = 6999.6665
$ws.Range("J67").Value = 6999.6665
$ws.Range("L67").Value = 6999.6665
$ws.Range("N67").Value = -8715.666499999999
# Row 98
$ws.Range("H98").Value = 2307.4285
$ws.Range("J98").Value = 3814.2856
$ws.Range("L98").Value = 3814.2856
$ws.Range("N98").Value = -6810.2856
# Row 106
$ws.Range("H106").Value = 63238.555
$ws.Range("I106").Value = 84605.414
$ws.Range("K106").Value = 84605.414
$ws.Range("M106").Value = -83974.414
# Row 112
$ws.Range("H112").Value = 1833.6774
$ws.Range("J112").Value = 1686.8889
$ws.Range("L112").Value = 5060.6667
$ws.Range("N112").Value = -7276.6667
# Row 122
$ws.Range("H122").Value = 2307.4285
$ws.Range("J122").Value = 3814.2856
$ws.Range("L122").Value = 11442.8568
$ws.Range("N122").Value = -16342.8568
# Row 131
$ws.Range("H131").Value = 1202.625
$ws.Range("I131").Value = 1202.625
$ws.Range("K131").Value = 3607.875
$ws.Range("M131").Value = 1432.125
# Row 137
$ws.Range("H137").Value = 4956.452
$ws.Range("I137").Value = 4689.5
$ws.Range("K137").Value = 14068.5
$ws.Range("M137").Value = -11518.5
# Row 138
$ws.Range("H138").Value = 3333.28
$ws.Range("J138").Value = 3724.4736
$ws.Range("L138").Value = 11173.4208
$ws.Range("N138").Value = -21453.4208
# Row 141
$ws.Range("H141").Value = 46881.6
$ws.Range("I141").Value = 51801.832
$ws.Range("K141").Value = 155405.496
$ws.Range("M141").Value = -150225.496

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6797533
$ws.Range("J32").Value = 33435812
$ws.Range("L32").Value = 33435812
$ws.Range("N32").Value = -33436386
# Row 74
$ws.Range("H74").Value = 1552.025
$ws.Range("J74").Value = 3655.6667
$ws.Range("L74").Value = 3655.6667
$ws.Range("N74").Value = -5403.6667
# Row 77
$ws.Range("H77").Value = 1552.025
$ws.Range("J77").Value = 3655.6667
$ws.Range("L77").Value = 18278.3335
$ws.Range("N77").Value = -27014.3335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 853.3182
$ws.Range("J80").Value = 870
$ws.Range("L80").Value = 870
$ws.Range("N80").Value = -2866
# Row 82
$ws.Range("H82").Value = 15955.818
# Row 83
$ws.Range("H83").Value = 853.3182
$ws.Range("J83").Value = 870
$ws.Range("L83").Value = 4350
$ws.Range("N83").Value = -14334
# Row 85
$ws.Range("H85").Value = 15955.818

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4478.5
$ws.Range("I31").Value = 2397.8
$ws.Range("J31").Value = 6212.4165
$ws.Range("K31").Value = 2397.8
$ws.Range("L31").Value = 6212.4165
$ws.Range("M31").Value = -2102.8
$ws.Range("N31").Value = -6802.4165
# Row 34
$ws.Range("H34").Value = 4478.5
$ws.Range("I34").Value = 2397.8
$ws.Range("J34").Value = 6212.4165
$ws.Range("K34").Value = 2397.8
$ws.Range("L34").Value = 6212.4165
$ws.Range("M34").Value = -2195.8
$ws.Range("N34").Value = -6616.4165
# Row 58
$ws.Range("H58").Value = 890.4
$ws.Range("I58").Value = 890.4
$ws.Range("K58").Value = 890.4
$ws.Range("M58").Value = -687.4
# Row 86
$ws.Range("H86").Value = 29416148
$ws.Range("I86").Value = 47622930
$ws.Range("K86").Value = 47622930
$ws.Range("M86").Value = -47621807
# Row 89
$ws.Range("H89").Value = 29416148
$ws.Range("I89").Value = 47622930
$ws.Range("K89").Value = 238114650
$ws.Range("M89").Value = -238109034
# Row 132
$ws.Range("H132").Value = 5755.25
$ws.Range("I132").Value = 5755.25
$ws.Range("K132").Value = 17265.75
$ws.Range("M132").Value = -14735.75
# Row 136
$ws.Range("H136").Value = 890.4
$ws.Range("I136").Value = 890.4
$ws.Range("K136").Value = 2671.2
$ws.Range("M136").Value = -121.1999999999998
# Row 137
$ws.Range("H137").Value = 38911.145
$ws.Range("J137").Value = 44992
$ws.Range("L137").Value = 44992
$ws.Range("N137").Value = -55192

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 184.38889
$ws.Range("J12").Value = 184.38889
$ws.Range("L12").Value = 553.1666700000001
$ws.Range("N12").Value = -899.1666700000001
# Row 139
$ws.Range("H139").Value = 49179.668
$ws.Range("I139").Value = 53672.26
$ws.Range("K139").Value = 161016.78
$ws.Range("M139").Value = -155876.78

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 11186
$ws.Range("I80").Value = 16461.143
$ws.Range("K80").Value = 16461.143
$ws.Range("M80").Value = -15463.143
# Row 83
$ws.Range("H83").Value = 11186
$ws.Range("I83").Value = 16461.143
$ws.Range("K83").Value = 82305.715
$ws.Range("M83").Value = -77313.715
# Row 102
$ws.Range("H102").Value = 2529.4443
$ws.Range("I102").Value = 2531.24
$ws.Range("J102").Value = 2507
$ws.Range("K102").Value = 2531.24
$ws.Range("L102").Value = 2507
$ws.Range("M102").Value = -909.2399999999998
$ws.Range("N102").Value = -5751

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 33372932
$ws.Range("I61").Value = 41673500
$ws.Range("K61").Value = 41673500
$ws.Range("M61").Value = -41673298
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = ""
$ws.Range("N63").Value = 0
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = ""
$ws.Range("N66").Value = 0
# Row 74
$ws.Range("H74").Value = 19000
$ws.Range("I74").Value = 19000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 19000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -18002
# Row 77
$ws.Range("H77").Value = 19000
$ws.Range("I77").Value = 19000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 57000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -52008
$ws.Range("N77").Value = ""
# Row 113
$ws.Range("H113").Value = 33372932
$ws.Range("I113").Value = 41673500
$ws.Range("K113").Value = 41673500
$ws.Range("M113").Value = -41671330
# Row 136
$ws.Range("H136").Value = 4016.4119
$ws.Range("I136").Value = 3254.889
$ws.Range("K136").Value = 9764.667000000001
$ws.Range("M136").Value = -7214.667000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3274.8125
$ws.Range("I122").Value = 2413.8572
$ws.Range("J122").Value = 3944.4443
$ws.Range("K122").Value = 7241.571599999999
$ws.Range("L122").Value = 11833.3329
$ws.Range("M122").Value = -4791.571599999999
$ws.Range("N122").Value = -16733.3329
# Row 132
$ws.Range("H132").Value = 17167.188
$ws.Range("I132").Value = 18191.5
$ws.Range("K132").Value = 54574.5
$ws.Range("M132").Value = -52044.5
# Row 136
$ws.Range("H136").Value = 8351
$ws.Range("J136").Value = 21248.5
$ws.Range("L136").Value = 63745.5
$ws.Range("N136").Value = -68845.5
